# Scheduled runner update: refresh cached Universalis market-price snapshots
# (currentAveragePrice / LevePrice / LeveProfit columns H:N) across all job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 731.3333
$ws.Range("I8").Value = 84.85714
$ws.Range("K8").Value = 254.57142
$ws.Range("M8").Value = -115.57142
$ws.Range("H86").Value = 2082.9333
$ws.Range("I86").Value = 1728.8
$ws.Range("K86").Value = 1728.8
$ws.Range("M86").Value = -605.8
$ws.Range("H89").Value = 2082.9333
$ws.Range("I89").Value = 1728.8
$ws.Range("K89").Value = 8644
$ws.Range("M89").Value = -3028
$ws.Range("H101").Value = 942
$ws.Range("J101").Value = 998.5
$ws.Range("L101").Value = 2995.5
$ws.Range("N101").Value = -6239.5
$ws.Range("H132").Value = 2297.3547
$ws.Range("I132").Value = 2297.3547
$ws.Range("K132").Value = 6892.0641
$ws.Range("M132").Value = -4362.0641
$ws.Range("H135").Value = 3339.7585
$ws.Range("I135").Value = 3650.6667
$ws.Range("J135").Value = 2523.625
$ws.Range("K135").Value = 32856.0003
$ws.Range("L135").Value = 22712.625
$ws.Range("M135").Value = -30321.0003
$ws.Range("N135").Value = -27782.625
$ws.Range("H138").Value = 5721.6177
$ws.Range("J138").Value = 5935.2593
$ws.Range("L138").Value = 17805.7779
$ws.Range("N138").Value = -28085.7779

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9523.200000000001
$ws.Range("I45").Value = 11722.546
$ws.Range("K45").Value = 11722.546
$ws.Range("M45").Value = -11345.546
$ws.Range("H61").Value = 4301.661
$ws.Range("J61").Value = 12012.25
$ws.Range("L61").Value = 12012.25
$ws.Range("N61").Value = -12436.25
$ws.Range("H74").Value = 5146.826
$ws.Range("I74").Value = 873.4666999999999
$ws.Range("K74").Value = 873.4666999999999
$ws.Range("M74").Value = 0.5333000000000538
$ws.Range("H77").Value = 5146.826
$ws.Range("I77").Value = 873.4666999999999
$ws.Range("K77").Value = 4367.3335
$ws.Range("M77").Value = 0.6665000000002692
$ws.Range("H122").Value = 1603.409
$ws.Range("I122").Value = 1408.1177
$ws.Range("K122").Value = 4224.3531
$ws.Range("M122").Value = -1774.3531
$ws.Range("H136").Value = 4301.661
$ws.Range("J136").Value = 12012.25
$ws.Range("L136").Value = 36036.75
$ws.Range("N136").Value = -41136.75

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 246332.67
$ws.Range("J70").Value = 246332.67
$ws.Range("L70").Value = 246332.67
$ws.Range("N70").Value = -246918.67
$ws.Range("H73").Value = 246332.67
$ws.Range("J73").Value = 246332.67
$ws.Range("L73").Value = 246332.67
$ws.Range("N73").Value = -248360.67
$ws.Range("H107").Value = 3526.516
$ws.Range("I107").Value = 3993.7
$ws.Range("J107").Value = 2677.0908
$ws.Range("K107").Value = 3993.7
$ws.Range("L107").Value = 2677.0908
$ws.Range("M107").Value = -2073.7
$ws.Range("N107").Value = -6517.0908

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1723.2727
$ws.Range("I94").Value = 4825.3335
$ws.Range("J94").Value = 560
$ws.Range("K94").Value = 4825.3335
$ws.Range("L94").Value = 560
$ws.Range("M94").Value = -4374.3335
$ws.Range("N94").Value = -1462
$ws.Range("H99").Value = 7663.5933
$ws.Range("I99").Value = 3590.25
$ws.Range("K99").Value = 3590.25
$ws.Range("M99").Value = -2092.25
$ws.Range("H126").Value = 7663.5933
$ws.Range("I126").Value = 3590.25
$ws.Range("K126").Value = 10770.75
$ws.Range("M126").Value = -8300.75
$ws.Range("H141").Value = 294580.8
$ws.Range("J141").Value = 307552.88
$ws.Range("L141").Value = 307552.88
$ws.Range("N141").Value = -317912.88

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1023.8
$ws.Range("I29").Value = 29
$ws.Range("J29").Value = 1687
$ws.Range("K29").Value = 87
$ws.Range("L29").Value = 5061
$ws.Range("M29").Value = 190
$ws.Range("N29").Value = -5615
$ws.Range("H46").Value = 9858478
$ws.Range("J46").Value = 1669.8
$ws.Range("L46").Value = 5009.4
$ws.Range("N46").Value = -5191.4
$ws.Range("H108").Value = 4517.5
$ws.Range("I108").Value = 4517.5
$ws.Range("K108").Value = 13552.5
$ws.Range("M108").Value = -10672.5
$ws.Range("H120").Value = 8553.888999999999
$ws.Range("I120").Value = 7830.8335
$ws.Range("K120").Value = 23492.5005
$ws.Range("M120").Value = -18654.5005

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7622.5483
$ws.Range("I70").Value = 6718.136
$ws.Range("J70").Value = 9833.333000000001
$ws.Range("K70").Value = 6718.136
$ws.Range("L70").Value = 9833.333000000001
$ws.Range("M70").Value = -6448.136
$ws.Range("N70").Value = -10373.333
$ws.Range("H73").Value = 7622.5483
$ws.Range("I73").Value = 6718.136
$ws.Range("J73").Value = 9833.333000000001
$ws.Range("K73").Value = 6718.136
$ws.Range("L73").Value = 9833.333000000001
$ws.Range("M73").Value = -5782.136
$ws.Range("N73").Value = -11705.333
$ws.Range("H80").Value = 1760.25
$ws.Range("I80").Value = 1697.4
$ws.Range("J80").Value = 1865
$ws.Range("K80").Value = 1697.4
$ws.Range("L80").Value = 1865
$ws.Range("M80").Value = -699.4000000000001
$ws.Range("N80").Value = -3861
$ws.Range("H83").Value = 1760.25
$ws.Range("I83").Value = 1697.4
$ws.Range("J83").Value = 1865
$ws.Range("K83").Value = 8487
$ws.Range("L83").Value = 9325
$ws.Range("M83").Value = -3495
$ws.Range("N83").Value = -19309
$ws.Range("H113").Value = 108047.9
$ws.Range("H122").Value = 1062.25
$ws.Range("I122").Value = 1099.9333
$ws.Range("J122").Value = 949.2
$ws.Range("K122").Value = 3299.7999
$ws.Range("L122").Value = 2847.6
$ws.Range("M122").Value = -849.7999
$ws.Range("N122").Value = -7747.6
$ws.Range("H132").Value = 3064.8286
$ws.Range("I132").Value = 2084.0667
$ws.Range("J132").Value = 8949.4
$ws.Range("K132").Value = 6252.2001
$ws.Range("L132").Value = 26848.2
$ws.Range("M132").Value = -3722.2001
$ws.Range("N132").Value = -31908.2

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5410.3076
$ws.Range("I40").Value = 5410.3076
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5410.3076
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5274.3076
$ws.Range("N40").Value = $null
$ws.Range("H51").Value = 42993.668
$ws.Range("J51").Value = 42993.668
$ws.Range("L51").Value = 42993.668
$ws.Range("N51").Value = -43949.668
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("H68").Value = 4749.5
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251
$ws.Range("H71").Value = 4749.5
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("H117").Value = 99999
$ws.Range("J117").Value = 99999
$ws.Range("L117").Value = 99999
$ws.Range("N117").Value = -109177

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8002.5
$ws.Range("J62").Value = 8002.5
$ws.Range("L62").Value = 8002.5
$ws.Range("N62").Value = -9250.5
$ws.Range("H64").Value = 59990
$ws.Range("I64").Value = 59989
$ws.Range("J64").Value = 59990.5
$ws.Range("K64").Value = 59989
$ws.Range("L64").Value = 59990.5
$ws.Range("M64").Value = -59741
$ws.Range("N64").Value = -60486.5
$ws.Range("H65").Value = 8002.5
$ws.Range("J65").Value = 8002.5
$ws.Range("L65").Value = 40012.5
$ws.Range("N65").Value = -46252.5
$ws.Range("H67").Value = 59990
$ws.Range("I67").Value = 59989
$ws.Range("J67").Value = 59990.5
$ws.Range("K67").Value = 59989
$ws.Range("L67").Value = 59990.5
$ws.Range("M67").Value = -59131
$ws.Range("N67").Value = -61706.5
$ws.Range("H75").Value = 23940
$ws.Range("I75").Value = 23940
$ws.Range("K75").Value = 23940
$ws.Range("M75").Value = -23004
$ws.Range("H78").Value = 23940
$ws.Range("I78").Value = 23940
$ws.Range("K78").Value = 71820
$ws.Range("M78").Value = -67140
$ws.Range("H94").Value = 51800
$ws.Range("J94").Value = 51800
$ws.Range("L94").Value = 51800
$ws.Range("N94").Value = -53602

